# PROGE_Time recording log.xlsx
# - Duplicate the empty "Nädal" template sheet, place the copy right before
#   it, rename the copy to "Nädal8" and fill in week 8's first log entry.
# - The new sheet becomes the active/selected tab; the previously active
#   sheet ("Nädal7") and the template sheet ("Nädal") lose their selection
#   highlighting and get their cursor moved elsewhere (as Excel naturally
#   does when the user clicks away to a different sheet/cell).

$wb = $excel.ActiveWorkbook

# 1. Duplicate the "Nädal" template sheet, inserting the copy immediately
#    before it (this keeps "Nädal" as the last / template sheet, matching
#    the tab order Nädal7, Nädal8, Nädal).
$template = $wb.Worksheets.Item("Nädal")
$template.Copy($template)

# The freshly-inserted copy sits right before the template, i.e. right
# after "Nädal7" -> position 8 (1-based) in the tab strip.
$newSheet = $wb.Worksheets.Item(8)
$newSheet.Name = "Nädal8"

# 2. Fill in the first time-log row (row 7) of the new sheet with week 8's
#    first entry.
$newSheet.Range("B7").Value = 43907
$newSheet.Range("C7").Value = 0.84166666666666667
$newSheet.Range("D7").Value = 0.90277777777777779
$newSheet.Range("G7").Value = "VL 37"
$newSheet.Range("H7").Value = "Master-Detail lõpetamine"

# 3. Move the cursor / selection around as it would naturally end up after
#    this editing session: land on E7 in the new sheet, move off Nädal7's
#    old selection, and leave the template sheet's cursor on G9.
$ws7 = $wb.Worksheets.Item("Nädal7")
$ws7.Range("D29").Select()

# Re-fetch the template sheet by name: after Copy() the old $template
# reference can be stale, so look it up fresh before selecting on it.
$templateAfter = $wb.Worksheets.Item("Nädal")
$templateAfter.Range("G9").Select()

$newSheet.Activate()
$newSheet.Range("E7").Select()
